# Flip the "Präsenz" (column C) boolean flag from FALSE to TRUE for the
# specific rows called out in the commit's updated data file.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Keyword Analysis")

$rows = @(2,10,16,23,31,32,33,67,72,77,89,98,146,168,169,170,176,179,182,186,188,189,190,194,197,198,206,250,256,264,289,311,312,313,338,341,350,366,470,474,477,488,514,515,530,531,562,569,570,571,610,614,615,618,621,674,675,689,725)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = $true
}
